$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top; everything currently in rows 1-28 shifts down to rows 4-31.
$ws.Range("A1:A3").EntireRow.Insert()

# Carry the number formats (styles) used by column A (date), column D (text) and
# columns F:G (text, left aligned) down into the freshly inserted rows, copying
# them from what is now row 4 (the original row 1) so no new style entries are
# created in styles.xml.
$ws.Range("A4").Copy()
$ws.Range("A1:A3").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D1:D3").PasteSpecial(-4122)
$ws.Range("F4:G4").Copy()
$ws.Range("F1:G3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 1: new debit movement - SPI COSTO OPER. CASH
$ws.Range("A1").Value = 41674
$ws.Range("B1").Value = "SPI COSTO OPER. CASH"
$ws.Range("C1").Value = "D"
$ws.Range("D1").Value = "0007868769"
$ws.Range("E1").Value = "CENTRO DE SERVIC. OPERAT. SS."
$ws.Range("F1").Value = "0.27  "
$ws.Range("G1").Value = "730.92"
$ws.Range("H1").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A1,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B1,"'', ''mo_tipo'' => ''",C1,"'', ''mo_documento'' => ''",D1,"'', ''mo_oficina'' => ''",E1,"'', ''mo_monto'' => ",F1,", ''mo_saldo'' => ",G1,", ''mo_fecha_crea'' => new \DateTime(''2014-02-01 00:00:01''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL),")'

# Row 2: new debit movement - IVA COBRADO
$ws.Range("A2").Value = 41674
$ws.Range("B2").Value = "IVA COBRADO"
$ws.Range("C2").Value = "D"
$ws.Range("D2").Value = "0007868769"
$ws.Range("E2").Value = "CENTRO DE SERVIC. OPERAT. SS."
$ws.Range("F2").Value = "0.03  "
$ws.Range("G2").Value = "731.19"
$ws.Range("H2").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A2,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B2,"'', ''mo_tipo'' => ''",C2,"'', ''mo_documento'' => ''",D2,"'', ''mo_oficina'' => ''",E2,"'', ''mo_monto'' => ",F2,", ''mo_saldo'' => ",G2,", ''mo_fecha_crea'' => new \DateTime(''2014-02-01 00:00:01''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL),")'

# Row 3: new credit movement - 50-SPI-CCU - MIN. ECONOMIA CUENTA
$ws.Range("A3").Value = 41674
$ws.Range("B3").Value = "50-SPI-CCU - MIN. ECONOMIA CUENTA"
$ws.Range("C3").Value = "C"
$ws.Range("D3").Value = "0007868747"
$ws.Range("E3").Value = "TENA"
$ws.Range("F3").Value = "724.27  "
$ws.Range("G3").Value = "731.22"
$ws.Range("H3").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A3,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B3,"'', ''mo_tipo'' => ''",C3,"'', ''mo_documento'' => ''",D3,"'', ''mo_oficina'' => ''",E3,"'', ''mo_monto'' => ",F3,", ''mo_saldo'' => ",G3,", ''mo_fecha_crea'' => new \DateTime(''2014-02-01 00:00:01''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL),")'

# Match the saved cursor/selection state: H1:H3 selected, active cell on H3.
$ws.Range("H1:H3").Select()
